$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.760.69'
$ws.Cells.Item(2, 5).Value = '  +1.68%  '

$ws.Cells.Item(3, 4).Value = '1.639.49'
$ws.Cells.Item(3, 5).Value = '  +2.27%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.996'
$ws.Cells.Item(4, 5).Value = '  -0.77%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '214.57'
$ws.Cells.Item(5, 5).Value = '  +0.77%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.504'
$ws.Cells.Item(6, 5).Value = '  +0.27%  '

$ws.Cells.Item(7, 5).Value = '  -0.87%  '

$ws.Cells.Item(8, 5).Value = '  +1.57%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.0612'
$ws.Cells.Item(9, 5).Value = '  +0.77%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.40'
$ws.Cells.Item(10, 5).Value = '  +2.33%  '

$ws.Cells.Item(11, 5).Value = '  +0.29%  '

$ws.Cells.Item(12, 5).Value = '  +1.75%  '

$ws.Cells.Item(13, 4).Value = '1.641.92'
$ws.Cells.Item(13, 5).Value = '  +2.37%  '

$ws.Cells.Item(14, 5).Value = '  +1.12%  '

$ws.Cells.Item(15, 5).Value = '  +2.03%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '64.43'
$ws.Cells.Item(16, 5).Value = '  +1.25%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '241.11'
$ws.Cells.Item(17, 5).Value = '  +6.22%  '

$ws.Cells.Item(18, 4).Value = '26.737.01'
$ws.Cells.Item(18, 5).Value = '  +1.53%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.85'
$ws.Cells.Item(19, 5).Value = '  +3.66%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.997'
$ws.Cells.Item(21, 5).Value = '  -0.74%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.36'
$ws.Cells.Item(22, 5).Value = '  +1.16%  '

$ws.Cells.Item(23, 5).Value = '  +2.34%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.21'
$ws.Cells.Item(24, 5).Value = '  +2.56%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '146.47'
$ws.Cells.Item(25, 5).Value = '  +0.74%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.997'
$ws.Cells.Item(26, 5).Value = '  -0.75%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '7.10'
$ws.Cells.Item(27, 5).Value = '  +2.31%  '

$ws.Cells.Item(28, 5).Value = '  +0.66%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '15.83'
$ws.Cells.Item(29, 5).Value = '  +2.70%  '

$ws.Cells.Item(30, 5).Value = '  +0.79%  '

$ws.Cells.Item(31, 5).Value = '  +0.18%  '

$ws.Cells.Item(32, 4).Value = '1.522.13'
$ws.Cells.Item(32, 5).Value = '  +5.14%  '

$ws.Cells.Item(33, 5).Value = '  +1.75%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.04'
$ws.Cells.Item(34, 5).Value = '  +2.44%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.55'
$ws.Cells.Item(35, 5).Value = '  +5.55%  '

$ws.Cells.Item(36, 5).Value = '  -0.52%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.574'
$ws.Cells.Item(37, 5).Value = '  +1.87%  '

$ws.Cells.Item(38, 5).Value = '  +1.46%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.854'
$ws.Cells.Item(39, 5).Value = '  +3.47%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.90'
$ws.Cells.Item(40, 5).Value = '  +1.16%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.997'
$ws.Cells.Item(41, 5).Value = '  -0.81%  '

$ws.Cells.Item(42, 5).Value = '  +2.27%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '63.78'
$ws.Cells.Item(43, 5).Value = '  +4.72%  '

$ws.Cells.Item(44, 4).Value = '1.777.37'
$ws.Cells.Item(44, 5).Value = '  +1.83%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.768'
$ws.Cells.Item(45, 5).Value = '  +1.06%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.907'
$ws.Cells.Item(46, 5).Value = '  -2.42%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '90.44'
$ws.Cells.Item(47, 5).Value = '  +3.29%  '

$ws.Cells.Item(48, 5).Value = '  +3.12%  '

$ws.Cells.Item(49, 5).Value = '  +0.23%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0974'
$ws.Cells.Item(50, 5).Value = '  +2.53%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.61'
$ws.Cells.Item(51, 5).Value = '  +3.12%  '
